# refactor: add db stock json and api update
#
# Adds two new stock-item rows ("Fresh Milk Green Field" / ml,
# "UHT Milk Indomilk" / ml) below the existing "Proof Blend Coffee" / gr
# row, reusing the existing column formats, and re-centers / re-fonts the
# unit ("B") column header style so the new "ml" unit cells line up with
# the existing "gr" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- update the shared "unit" column style (column B) -----------------
# Previously Arial 10 with default alignment; now Calibri 11, centered
# horizontally, bottom-aligned vertically (matches column A's font).
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B1").VerticalAlignment   = -4107   # xlBottom
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 11

# --- row 2: Fresh Milk Green Field / ml --------------------------------
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)           # xlPasteFormats
$ws.Range("A2").Value = "Fresh Milk Green Field"

$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)           # xlPasteFormats
$ws.Range("B2").Value = "ml"

# --- row 3: UHT Milk Indomilk / ml -------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122)           # xlPasteFormats
$ws.Range("A3").Value = "UHT Milk Indomilk"

$ws.Range("B1").Copy()
$ws.Range("B3").PasteSpecial(-4122)           # xlPasteFormats
$ws.Range("B3").Value = "ml"
